$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in columns A and B for rows 2-5
$ws.Range("A2").Value = 44876
$ws.Range("B2").Value = 44906

$ws.Range("A3").Value = 44876
$ws.Range("B3").Value = 44906

$ws.Range("A4").Value = 44876
$ws.Range("B4").Value = 44906

$ws.Range("A5").Value = 44876
$ws.Range("B5").Value = 44906

# Update the selection to A2:E5 with active cell A2
$ws.Range("A2:E5").Select()
